# Add "(Browser on OS)" context suffixes to the Tool/Screen-reader header
# row (row 3) so the column headers are self-describing for screen reader
# users, per the commit message:
#   "Added additional information in column headers for making them
#    screen reader friendly"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Windows / Chrome block (B:G)
$ws.Range("B3").Value = "builtin inspector (Chrome on Windows)"
$ws.Range("C3").Value = "accessibility insights (Chrome on Windows)"
$ws.Range("D3").Value = "JAWS 2021 (Chrome on Windows)"
$ws.Range("E3").Value = "spoken by JAWS (Chrome on Windows)"
$ws.Range("F3").Value = "NVDA 2021 (Chrome on Windows)"
$ws.Range("G3").Value = "spoken by NVDA (Chrome on Windows)"

# Windows / Firefox block (H:M)
$ws.Range("H3").Value = "inspector (Firefox on Windows)"
$ws.Range("I3").Value = "accessibility insights (Firefox on Windows)"
$ws.Range("J3").Value = "JAWS 2021 (Firefox on Windows)"
$ws.Range("K3").Value = "spoken by JAWS (Firefox on Windows)"
$ws.Range("L3").Value = "NVDA 2020 (Firefox on Windows)"
$ws.Range("M3").Value = "spoken by NVDA (Firefox on Windows)"

# Mac / Safari block (N:O)
$ws.Range("N3").Value = "developer console (Safari on Mac)"
$ws.Range("O3").Value = "spoken by VoiceOver (Safari on Mac)"

# Mac / Chrome block (P:Q)
$ws.Range("P3").Value = "developer console (Chrome on Mac)"
$ws.Range("Q3").Value = "spoken by VoiceOver (Chrome on Mac)"

# Mac / Firefox block (R:S)
$ws.Range("R3").Value = "developer console (Firefox on Mac)"
$ws.Range("S3").Value = "spoken by VoiceOver (Firefox on Mac)"

# iOS / Safari (T)
$ws.Range("T3").Value = "spoken by VoiceOver (Safari on iOS)"

# Android / Chrome (U)
$ws.Range("U3").Value = "spoken by TalkBack (Chrome on Android)"

# Update the saved view/selection state to match the author's working
# position at the time of the edit (scrolled down toward the pagebreak
# row, bottom-right pane selection on D32). The frozen-pane split
# (column A / rows 1-3) itself is left untouched.
$ws.Range("D32").Select()
